$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Products")

# Delete rows 3 to 6 (only row 2 of data should remain)
$ws.Rows("3:6").Delete() | Out-Null

# Update header row: insert new columns H (barcodePath) and I (barcodeType)
# before the existing createdAt/updatedAt columns.
$ws.Range("H1").Value = "barcodePath"
$ws.Range("I1").Value = "barcodeType"
$ws.Range("J1").Value = "createdAt"
$ws.Range("K1").Value = "updatedAt"

# Update the single remaining data row (row 2) with the new values
$ws.Range("A2").Value = "dba28441-3044-46f2-a914-204bf9fd0954"
$ws.Range("B2").Value = "Test"
$ws.Range("C2").Value = "Buah"
$ws.Range("D2").Value = 123
$ws.Range("E2").Value = 123
$ws.Range("F2").Value = "T"
$ws.Range("G2").Value = "D:\Project\Developments\inventory\database\images\qr_T.png"
$ws.Range("H2").Value = "D:\Project\Developments\inventory\database\images\barcode_T.png"
$ws.Range("I2").Value = "code128"
$ws.Range("J2").Value = "2025-07-01T18:13:10.810Z"
$ws.Range("K2").Value = "2025-07-01T18:22:09.531Z"
